$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1176.3334
$ws.Range("I129").Value = 528.3333
$ws.Range("J129").Value = 1407.762
$ws.Range("K129").Value = 1584.9999
$ws.Range("L129").Value = 4223.286
$ws.Range("M129").Value = 3415.0001
$ws.Range("N129").Value = -14223.286

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 610.7308
$ws.Range("I97").Value = 741
$ws.Range("J97").Value = 176.5
$ws.Range("K97").Value = 741
$ws.Range("L97").Value = 176.5
$ws.Range("M97").Value = -245
$ws.Range("N97").Value = -1168.5

# ARM row 102
$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1217.8334
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 1226.75
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 1226.75
$ws.Range("M86").Value = -77
$ws.Range("N86").Value = -3472.75

# BSM row 89
$ws.Range("H89").Value = 1217.8334
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 1226.75
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 6133.75
$ws.Range("M89").Value = -384
$ws.Range("N89").Value = -17365.75

# BSM row 94
$ws.Range("H94").Value = 1209.8334
$ws.Range("I94").Value = 1234.2222
$ws.Range("J94").Value = 1136.6666
$ws.Range("K94").Value = 1234.2222
$ws.Range("L94").Value = 1136.6666
$ws.Range("M94").Value = -783.2221999999999
$ws.Range("N94").Value = -2038.6666

# BSM row 99
$ws.Range("H99").Value = 2760
$ws.Range("I99").Value = 1393.3334
$ws.Range("J99").Value = 3443.3333
$ws.Range("K99").Value = 1393.3334
$ws.Range("L99").Value = 3443.3333
$ws.Range("M99").Value = 104.6666
$ws.Range("N99").Value = -6439.3333

# BSM row 105
$ws.Range("H105").Value = 2480
$ws.Range("J105").Value = 3500
$ws.Range("L105").Value = 3500
$ws.Range("N105").Value = -6994

# BSM row 134
$ws.Range("H134").Value = 4054.4119
$ws.Range("I134").Value = 1619.6666
$ws.Range("J134").Value = 11967.333
$ws.Range("K134").Value = 4858.9998
$ws.Range("L134").Value = 35901.999
$ws.Range("M134").Value = -2323.9998
$ws.Range("N134").Value = -40971.999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4506767.5
$ws.Range("I31").Value = 1767.6578
$ws.Range("J31").Value = 9262045
$ws.Range("K31").Value = 1767.6578
$ws.Range("L31").Value = 9262045
$ws.Range("M31").Value = -1472.6578
$ws.Range("N31").Value = -9262635

# CRP row 34
$ws.Range("H34").Value = 4506767.5
$ws.Range("I34").Value = 1767.6578
$ws.Range("J34").Value = 9262045
$ws.Range("K34").Value = 1767.6578
$ws.Range("L34").Value = 9262045
$ws.Range("M34").Value = -1565.6578
$ws.Range("N34").Value = -9262449

# CRP row 62
$ws.Range("H62").Value = 6388.4
$ws.Range("I62").Value = 2265.5
$ws.Range("J62").Value = 22880
$ws.Range("K62").Value = 2265.5
$ws.Range("L62").Value = 22880
$ws.Range("M62").Value = -1641.5
$ws.Range("N62").Value = -24128

# CRP row 65
$ws.Range("H65").Value = 6388.4
$ws.Range("I65").Value = 2265.5
$ws.Range("J65").Value = 22880
$ws.Range("K65").Value = 11327.5
$ws.Range("L65").Value = 114400
$ws.Range("M65").Value = -8207.5
$ws.Range("N65").Value = -120640

# CRP row 105
$ws.Range("H105").Value = 1085.6666
$ws.Range("I105").Value = 1085.6666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1085.6666
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = 661.3334

# CRP row 134
$ws.Range("H134").Value = 1602.8043
$ws.Range("I134").Value = 1195.6957
$ws.Range("J134").Value = 2009.9131
$ws.Range("K134").Value = 3587.0871
$ws.Range("L134").Value = 6029.7393
$ws.Range("M134").Value = -1052.0871
$ws.Range("N134").Value = -11099.7393

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2715.3333
$ws.Range("I80").Value = 2680
$ws.Range("J80").Value = 2750.6667
$ws.Range("K80").Value = 2680
$ws.Range("L80").Value = 2750.6667
$ws.Range("M80").Value = -1682
$ws.Range("N80").Value = -4746.6667

# GSM row 83
$ws.Range("H83").Value = 2715.3333
$ws.Range("I83").Value = 2680
$ws.Range("J83").Value = 2750.6667
$ws.Range("K83").Value = 13400
$ws.Range("L83").Value = 13753.3335
$ws.Range("M83").Value = -8408
$ws.Range("N83").Value = -23737.3335

# GSM row 97
$ws.Range("H97").Value = 926.1923
$ws.Range("I97").Value = 967.34784
$ws.Range("J97").Value = 610.6667
$ws.Range("K97").Value = 967.34784
$ws.Range("L97").Value = 610.6667
$ws.Range("M97").Value = -471.34784
$ws.Range("N97").Value = -1602.6667

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2259.0908
$ws.Range("I68").Value = 1678.5714
$ws.Range("J68").Value = 3275
$ws.Range("K68").Value = 1678.5714
$ws.Range("L68").Value = 3275
$ws.Range("M68").Value = -929.5714
$ws.Range("N68").Value = -4773

# LTW row 71
$ws.Range("H71").Value = 2259.0908
$ws.Range("I71").Value = 1678.5714
$ws.Range("J71").Value = 3275
$ws.Range("K71").Value = 8392.857
$ws.Range("L71").Value = 16375
$ws.Range("M71").Value = -4648.857
$ws.Range("N71").Value = -23863

# LTW row 82
$ws.Range("H82").Value = 1700.8572
$ws.Range("I82").Value = 1273.1428
$ws.Range("J82").Value = 2128.5715
$ws.Range("K82").Value = 1273.1428
$ws.Range("L82").Value = 2128.5715
$ws.Range("M82").Value = -912.1428000000001
$ws.Range("N82").Value = -2850.5715

# LTW row 85
$ws.Range("H85").Value = 1700.8572
$ws.Range("I85").Value = 1273.1428
$ws.Range("J85").Value = 2128.5715
$ws.Range("K85").Value = 1273.1428
$ws.Range("L85").Value = 2128.5715
$ws.Range("M85").Value = -25.14280000000008
$ws.Range("N85").Value = -4624.5715

# LTW row 93
$ws.Range("H93").Value = 2219.8572
$ws.Range("I93").Value = 1925
$ws.Range("J93").Value = 2514.7144
$ws.Range("K93").Value = 1925
$ws.Range("L93").Value = 2514.7144
$ws.Range("M93").Value = -677
$ws.Range("N93").Value = -5010.7144

# LTW row 100
$ws.Range("H100").Value = 2730.4
$ws.Range("I100").Value = 1350
$ws.Range("J100").Value = 3650.6667
$ws.Range("K100").Value = 1350
$ws.Range("L100").Value = 3650.6667
$ws.Range("M100").Value = -809
$ws.Range("N100").Value = -4732.6667

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4324.4
$ws.Range("I62").Value = 3183.5
$ws.Range("K62").Value = 3183.5
$ws.Range("M62").Value = -2559.5

# WVR row 65
$ws.Range("H65").Value = 4324.4
$ws.Range("I65").Value = 3183.5
$ws.Range("K65").Value = 15917.5
$ws.Range("M65").Value = -12797.5

# WVR row 81
$ws.Range("H81").Value = 1170.1562
$ws.Range("I81").Value = 1433.9375
$ws.Range("J81").Value = 906.375
$ws.Range("K81").Value = 2867.875
$ws.Range("L81").Value = 1812.75
$ws.Range("M81").Value = -1806.875
$ws.Range("N81").Value = -3934.75

# WVR row 84
$ws.Range("H84").Value = 1170.1562
$ws.Range("I84").Value = 1433.9375
$ws.Range("J84").Value = 906.375
$ws.Range("K84").Value = 14339.375
$ws.Range("L84").Value = 9063.75
$ws.Range("M84").Value = -9035.375
$ws.Range("N84").Value = -19671.75

# WVR row 96
$ws.Range("H96").Value = 1522
$ws.Range("I96").Value = 1047.1428
$ws.Range("J96").Value = 2186.8
$ws.Range("K96").Value = 1047.1428
$ws.Range("L96").Value = 2186.8
$ws.Range("M96").Value = 325.8571999999999
$ws.Range("N96").Value = -4932.8
